$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row originally had a duplicate "t1 [s]" label in both B1 and C1.
# Fix the labelling so each measurement column gets its own sequential label:
# B1 stays "t1 [s]"; C1..K1 shift up to the existing "t2 [s]" .. "t10 [s]"
# labels, and L1 becomes a brand new "t11 [s]" label.
$ws.Range("C1").Value2 = "t2 [s]"
$ws.Range("D1").Value2 = "t3 [s]"
$ws.Range("E1").Value2 = "t4 [s]"
$ws.Range("F1").Value2 = "t5 [s]"
$ws.Range("G1").Value2 = "t6 [s]"
$ws.Range("H1").Value2 = "t7 [s]"
$ws.Range("I1").Value2 = "t8 [s]"
$ws.Range("J1").Value2 = "t9 [s]"
$ws.Range("K1").Value2 = "t10 [s]"
$ws.Range("L1").Value2 = "t11 [s]"

# Match the existing header style: the leading "t" keeps the default run
# formatting, the number is shown smaller (8pt) and the " [s]" suffix uses
# the normal 11pt size, same as the other header cells in the row.
$ws.Range("L1").Characters(2, 2).Font.Size = 8
$ws.Range("L1").Characters(4, 4).Font.Size = 11

# Leave the selection on the newly-labelled last column.
[void]$ws.Range("L1").Select()
